$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "linha" values for alternating rows to "Linha 2"
$ws.Range("B3").Value = "Linha 2"
$ws.Range("B5").Value = "Linha 2"
$ws.Range("B7").Value = "Linha 2"
$ws.Range("B9").Value = "Linha 2"
$ws.Range("B11").Value = "Linha 2"
$ws.Range("B13").Value = "Linha 2"

# Update the active selection to J16
$ws.Range("J16").Select()
